$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4488961696624756
$ws.Range("B1").Value = 3.523543119430542
$ws.Range("C1").Value = 4.117266654968262
$ws.Range("D1").Value = 1.359105706214905
$ws.Range("E1").Value = 0.899681568145752
